$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while keeping it as Text (matches original inlineStr
# cells) and without leaving a residual non-default cell style behind.
function Set-TextValue {
    param($Sheet, $Address, $Text)
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '29.389.39'
Set-TextValue $ws 'E2' '  -0.33%  '
Set-TextValue $ws 'D3' '1.847.39'
Set-TextValue $ws 'E3' '  -0.12%  '
Set-TextValue $ws 'E4' '  -0.09%  '
Set-TextValue $ws 'D5' '240.28'
Set-TextValue $ws 'E5' '  -0.67%  '
Set-TextValue $ws 'D6' '0.6321'
Set-TextValue $ws 'E6' '  +0.88%  '
Set-TextValue $ws 'D7' '0.9999'
Set-TextValue $ws 'D8' '0.07561'
Set-TextValue $ws 'E8' '  +0.32%  '
Set-TextValue $ws 'D9' '0.2962'
Set-TextValue $ws 'E9' '  -0.33%  '
Set-TextValue $ws 'D10' '24.47'
Set-TextValue $ws 'E10' '  +0.84%  '
Set-TextValue $ws 'D11' '0.07713'
Set-TextValue $ws 'E11' '  +0.46%  '
Set-TextValue $ws 'D12' '1.850.05'
Set-TextValue $ws 'E12' '  -0.56%  '
Set-TextValue $ws 'D13' '4.996'
Set-TextValue $ws 'E13' '  -0.29%  '
Set-TextValue $ws 'D14' '0.6854'
Set-TextValue $ws 'E14' '  +0.07%  '
Set-TextValue $ws 'D15' '0.00001011'
Set-TextValue $ws 'E15' '  +3.96%  '
Set-TextValue $ws 'D16' '83.04'
Set-TextValue $ws 'E16' '  -0.92%  '
Set-TextValue $ws 'D17' '6.157'
Set-TextValue $ws 'E17' '  -1.17%  '
Set-TextValue $ws 'D18' '29.421.07'
Set-TextValue $ws 'E18' '  -0.39%  '
Set-TextValue $ws 'D19' '229.84'
Set-TextValue $ws 'E19' '  -1.92%  '
Set-TextValue $ws 'D20' '12.45'
Set-TextValue $ws 'E20' '  -0.33%  '
Set-TextValue $ws 'D21' '0.9996'
Set-TextValue $ws 'E21' '  -0.07%  '
Set-TextValue $ws 'D22' '7.558'
Set-TextValue $ws 'E23' '  +0.00%  '
Set-TextValue $ws 'D24' '156.85'
Set-TextValue $ws 'E24' '  +0.67%  '
Set-TextValue $ws 'D25' '0.1400'
Set-TextValue $ws 'E25' '  +0.89%  '
Set-TextValue $ws 'D26' '8.388'
Set-TextValue $ws 'E26' '  -0.42%  '
Set-TextValue $ws 'E27' '  -0.35%  '
Set-TextValue $ws 'E28' '  -0.98%  '
Set-TextValue $ws 'D29' '0.05734'
Set-TextValue $ws 'E29' '  -1.73%  '
Set-TextValue $ws 'E30' '  -0.69%  '
Set-TextValue $ws 'E31' '  +0.61%  '
Set-TextValue $ws 'D32' '4.028'
Set-TextValue $ws 'E32' '  -0.05%  '
Set-TextValue $ws 'E33' '  -2.29%  '
Set-TextValue $ws 'E34' '  -1.17%  '
Set-TextValue $ws 'D35' '0.7162'
Set-TextValue $ws 'D36' '2.590'
Set-TextValue $ws 'E36' '  +0.05%  '
Set-TextValue $ws 'D37' '1.253.77'
Set-TextValue $ws 'E37' '  +1.45%  '
Set-TextValue $ws 'E38' '  +2.49%  '
Set-TextValue $ws 'D39' '2.780'
Set-TextValue $ws 'E39' '  -0.60%  '
Set-TextValue $ws 'D40' '0.9100'
Set-TextValue $ws 'E40' '  -0.32%  '
Set-TextValue $ws 'D41' '6.176'
Set-TextValue $ws 'E41' '  +0.72%  '
Set-TextValue $ws 'E42' '  +0.02%  '
Set-TextValue $ws 'D43' '1.995.35'
Set-TextValue $ws 'E43' '  -2.40%  '
Set-TextValue $ws 'D44' '101.67'
Set-TextValue $ws 'E44' '  -0.69%  '
Set-TextValue $ws 'D45' '66.24'
Set-TextValue $ws 'E45' '  -1.57%  '
Set-TextValue $ws 'D46' '7.067'
Set-TextValue $ws 'E46' '  -2.96%  '
Set-TextValue $ws 'B47' 'BabyDogeCoin'
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws 'D47' '0.00000000118'
Set-TextValue $ws 'E47' '  +0.18%  '
Set-TextValue $ws 'B48' 'TheSandbox'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws 'D48' '0.4028'
Set-TextValue $ws 'E48' '  -0.06%  '
Set-TextValue $ws 'B49' 'EnergySwap'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 'D49' '9.126'
Set-TextValue $ws 'E49' '  -0.08%  '
Set-TextValue $ws 'D50' '1.703'
Set-TextValue $ws 'E50' '  +0.75%  '
Set-TextValue $ws 'E51' '  +1.04%  '
